$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 00:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 845822
$ws.Range("C4").Value = 27078
$ws.Range("D4").Value = 83917
$ws.Range("E4").Value = 714426
$ws.Range("G4").Value = 2161
$ws.Range("H4").Value = 47479

# Row 50 - Colombia
$ws.Range("B50").Value = 4356
$ws.Range("C50").Value = 207
$ws.Range("D50").Value = 870
$ws.Range("E50").Value = 3280
$ws.Range("G50").Value = 10
$ws.Range("H50").Value = 206

# Row 87 - Costa de Marfil
$ws.Range("B87").Value = 952
$ws.Range("C87").Value = 36
$ws.Range("D87").Value = 310
$ws.Range("E87").Value = 628
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 14
